# MAI_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer note from 2021-05-05 to 2021-05-06
#  - refresh the Weight / Percent Change figures for the 6 holdings rows (rows 2-7)
#
# The worksheet ships with cell protection enabled, so it must be unprotected
# before the cells can be written to, then re-protected to restore the sheet's
# original state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# --- Disclaimer note (row 10) -------------------------------------------------
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."

# --- Holdings table (columns D = Weight, E = Percent Change) -----------------
$ws.Range("D2").Value = 0.4805400013825932
$ws.Range("E2").Value = -0.0007809449433814342

$ws.Range("D3").Value = 0.3393500824618059
$ws.Range("E3").Value = 0.01040385888584128

$ws.Range("D4").Value = 0.09519153852990846
$ws.Range("E4").Value = 0.01053013798111846

$ws.Range("D5").Value = 0.05388409919118301
$ws.Range("E5").Value = -0.0003436426116839586

$ws.Range("D6").Value = 0.03103427843450953
$ws.Range("E6").Value = 0.01223150357995206

$ws.Range("E7").Value = 0.004518734137212599

# Restore worksheet protection.
$ws.Protect()
